$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F
$ws.Range("F1").Value = "Operant Start (seconds)"
$ws.Range("F1").Font.Bold = $true

# Data values for F2:F12 (video behavior offset, in seconds)
$values = @(-1, -0.5, -0.5, -0.5, -1, -1, -1, -0.5, -0.5, -0.5, -0.5)
$ws.Range("B2").Copy()
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 6)
    $cell.PasteSpecial(-4122)
    $cell.Value = $values[$i]
}

$ws.Range("G17").Select()
